$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.938207507133484
$ws.Range("B1").Value = 2.187227487564087
$ws.Range("C1").Value = 2.13093376159668
$ws.Range("D1").Value = 2.549808025360107
$ws.Range("E1").Value = 2.06237268447876
